$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the backlog item text that was added to column A / B (rows 4-23)
# ---------------------------------------------------------------------------
$ws.Range("A4").Value  = "Select games that I own/ enjoy"
$ws.Range("B4").Value  = "UI for selection, searchable, possibily sortable by genres"
$ws.Range("B5").Value  = "Game object, that holds game information"
$ws.Range("B6").Value  = "Genre object, which holds list of similar games"
$ws.Range("B7").Value  = "No database as of now so prexisting local list of games"
$ws.Range("B8").Value  = "User Validation(might not need since we are the ones uploading games)"
$ws.Range("B9").Value  = "Testing model/ viewmodel"
$ws.Range("A10").Value = "Find recommended/ suggested games"
$ws.Range("B10").Value = "UI for recommendation engine, swipe left or right"
$ws.Range("B11").Value = "UI element for each recommended game"
$ws.Range("B12").Value = "Functionality to save, liked/ dislikes games and genres, based on user's choices"
$ws.Range("B13").Value = "Testing model/ viewmodel"
$ws.Range("A14").Value = "Access and view game library"
$ws.Range("B14").Value = "UI for game library, show user's owned games"
$ws.Range("B15").Value = "UI to display details from game, use image here"
$ws.Range("B16").Value = "UI for ability to select new game user owns"
$ws.Range("B17").Value = "Library object, which holds maybe a user and then pull the list of games related to that user"
$ws.Range("B18").Value = "Testing model/ viewmodel"
$ws.Range("A19").Value = "Create a profile"
$ws.Range("B19").Value = "UI for creating a profile, needs to navigate to UI selection of games"
$ws.Range("B20").Value = "UI for basic profile page"
$ws.Range("B21").Value = "UI allowing user to edit preferences, undo likes, dislikes"
$ws.Range("B22").Value = "User/ Profile object, holds user information, list of own games, likes and dislikes"
$ws.Range("B23").Value = "Testing model/ viewmodel"

Write-Output "text filled"

# ---------------------------------------------------------------------------
# 2. Insert a new blank backlog row right above the old totals row (row 27),
#    pushing the totals row down to row 28 and extending the SUM() ranges so
#    they include the freshly inserted row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(27).Insert()

$ws.Cells.Item(28, 3).Formula = "=SUM(C3:C27)"
$ws.Range("D28:G28").Formula = "=SUM(D3:D27)"

Write-Output "row inserted, totals updated"

# ---------------------------------------------------------------------------
# 3. The "Estimate Totals" label used to live next to the SUM() row (B27 ->
#    after the insert above it is B28). It now moves to its own row (B30),
#    two rows further down, with a genuinely blank spacer row (29) between
#    the totals row and the label. Grab B28's current formatting (still the
#    original bold/right-aligned "Estimate Totals" style) and stamp it onto
#    B30 first, then clear/restyle B28 to match a plain data row.
# ---------------------------------------------------------------------------
$ws.Range("B28:B28").Copy() | Out-Null
$ws.Cells.Item(30, 2).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Cells.Item(30, 2).Value = "Estimate Totals"

$ws.Cells.Item(28, 2).ClearContents()
$ws.Range("A26:A26").Copy() | Out-Null
$ws.Cells.Item(28, 2).PasteSpecial(-4122) | Out-Null

# Row 29 is a blank spacer row; only column B carries the plain row style,
# matching the other "background" cells in column A/B.
$ws.Range("A26:A26").Copy() | Out-Null
$ws.Cells.Item(29, 2).PasteSpecial(-4122) | Out-Null

Write-Output "B28 cleared, spacer row 29 + totals label row 30 written"

# ---------------------------------------------------------------------------
# 4. The burndown chart's series referenced the old totals row (C27:G27);
#    repoint it at the new totals row (C28:G28).
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(,,Sheet1!`$C`$28:`$G`$28,1)"

Write-Output "chart series repointed"


